# Generate Report for Handoff
# Updates status text from "In Translation" to "Ready for handoff" and
# refreshes the associated timestamp columns across the Overview, zh-cn and
# de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 06:38:23"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 06:38:18"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 06:38:23"

# Widen the status columns so they keep fitting the new, longer text
# (mirrors the column width change seen in the diff, which grows each
# status column from ~13.41 to ~17.22 characters).
#
# Note: the host's ColumnWidth setter quantizes to 1/6-character steps
# (stored = round((input + 5/6) * 6) / 6), so 16.33333333333333 is the
# input value that lands on 17.16666666666667 -- the closest reachable
# step to the recorded 17.2159881591797 target width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.33333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33333333333333
